$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the verbose "Round" column (F) values with their abbreviated forms.
$ws.Range("F2").Value = "F"
$ws.Range("F3").Value = "SF"
$ws.Range("F4").Value = "F"
$ws.Range("F5").Value = "R16"
$ws.Range("F6").Value = "QF"
$ws.Range("F7").Value = "F"
$ws.Range("F8").Value = "F"

# Move the active selection to F9, matching the saved view state.
$ws.Range("F9").Select()
